# bioSample_2651.xlsx edit
# Commit message: "fixing strain number and entering marker info for off by one errors"
#
# Two related fixes applied to the "90minuteinduction" block for replicate 3
# (rows 11-21) and replicate 4 (rows 33-43):
#   1. The "strain" (column E) values in these rows were off-by-one versus the
#      genotype (column F) they were paired with, so each strain id is
#      decremented by one position (TDY2194->TDY2193, TDY2195->TDY2194, ...).
#      Row 18 / row 40 (bioSampleNumber 17 / 39) is the one exception: its
#      naive decrement would have produced a duplicate/typo id, so those two
#      rows are set explicitly (row 18 keeps the raw "TDY22010" artefact,
#      row 40 is hand corrected to "TDY2200").
#   2. A "NAT" marker (column J / marker_1) is recorded for every one of
#      those corrected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replicate 3 block (rows 11-21) ---------------------------------------
$ws.Range("E11").Value = "TDY2193"
$ws.Range("E12").Value = "TDY2194"
$ws.Range("E13").Value = "TDY2195"
$ws.Range("E14").Value = "TDY2196"
$ws.Range("E15").Value = "TDY2197"
$ws.Range("E16").Value = "TDY2198"
$ws.Range("E17").Value = "TDY2199"
$ws.Range("E18").Value = "TDY22010"
$ws.Range("E19").Value = "TDY2201"
$ws.Range("E20").Value = "TDY2202"
$ws.Range("E21").Value = "TDY2203"

$ws.Range("J11").Value = "NAT"
$ws.Range("J12").Value = "NAT"
$ws.Range("J13").Value = "NAT"
$ws.Range("J14").Value = "NAT"
$ws.Range("J15").Value = "NAT"
$ws.Range("J16").Value = "NAT"
$ws.Range("J17").Value = "NAT"
$ws.Range("J18").Value = "NAT"
$ws.Range("J19").Value = "NAT"
$ws.Range("J20").Value = "NAT"
$ws.Range("J21").Value = "NAT"

# --- Replicate 4 block (rows 33-43) ---------------------------------------
$ws.Range("E33").Value = "TDY2193"
$ws.Range("E34").Value = "TDY2194"
$ws.Range("E35").Value = "TDY2195"
$ws.Range("E36").Value = "TDY2196"
$ws.Range("E37").Value = "TDY2197"
$ws.Range("E38").Value = "TDY2198"
$ws.Range("E39").Value = "TDY2199"
$ws.Range("E40").Value = "TDY2200"
$ws.Range("E41").Value = "TDY2201"
$ws.Range("E42").Value = "TDY2202"
$ws.Range("E43").Value = "TDY2203"

$ws.Range("J33").Value = "NAT"
$ws.Range("J34").Value = "NAT"
$ws.Range("J35").Value = "NAT"
$ws.Range("J36").Value = "NAT"
$ws.Range("J37").Value = "NAT"
$ws.Range("J38").Value = "NAT"
$ws.Range("J39").Value = "NAT"
$ws.Range("J40").Value = "NAT"
$ws.Range("J41").Value = "NAT"
$ws.Range("J42").Value = "NAT"
$ws.Range("J43").Value = "NAT"

# Match the author's final cursor position recorded in the saved view state.
[void]$ws.Range("G12").Select()
